$wb = $excel.ActiveWorkbook

# Add the new "squrs" worksheet after the existing "wallet" sheet
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "squrs"

# Headers
$ws2.Range("A1").Value = "x"
$ws2.Range("B1").Value = "sqr(x)"

# Fill x values 1..18 and the sqr(x) formulas
for ($i = 1; $i -le 18; $i++) {
    $row = $i + 1
    $ws2.Cells.Item($row, 1).Value = $i
    $ws2.Cells.Item($row, 2).Formula = "=A$row*A$row"
}

# Selection on new sheet: B2
$ws2.Range("B2").Select()

# Original sheet selection becomes just A1 (was A:A) and it's no longer the active tab
$ws1.Range("A1").Select()

$ws2.Activate()
